$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 12 ("Metabuli") values with new Metabuli version results
$ws.Range("B12").Value = 0.87
$ws.Range("C12").Value = 0.67
$ws.Range("D12").Value = 0.98
$ws.Range("E12").Value = 0.65
$ws.Range("F12").Value = 0.54
$ws.Range("G12").Value = 0.85
$ws.Range("H12").Value = 0.71
$ws.Range("I12").Value = 0.76
$ws.Range("J12").Value = 0.85
$ws.Range("K12").Value = 0.78
$ws.Range("M12").Value = 0.91
$ws.Range("N12").Value = 0.83
$ws.Range("P12").Value = 0.95
